# Updated symbol list on Fri Dec 30 11:48:35 UTC 2022 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) stores numeric-looking values as plain text in the
# source workbook (t="inlineStr"). Writing a numeric-looking string via
# .Value normally gets auto-converted to a real number by Excel, which
# would change the stored cell type. Forcing the NumberFormat to Text
# ("@") before the write keeps it textual; ClearFormats() afterwards drops
# the now-unneeded explicit style again so the cell's formatting/style
# index is left exactly as it was originally (no style attribute).
function Set-TextCell($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

# Simple price (column D) updates
Set-TextCell "D2"  "243.53"
Set-TextCell "D3"  "25.20"
Set-TextCell "D4"  "5.164"
Set-TextCell "D5"  "0.05733"
Set-TextCell "D6"  "6.490"
Set-TextCell "D7"  "3.110"
Set-TextCell "D8"  "0.8091"
Set-TextCell "D9"  "0.8480"
Set-TextCell "D10" "0.1337"
Set-TextCell "D11" "0.06948"
Set-TextCell "D12" "0.02830"
Set-TextCell "D13" "0.09372"

# Row 15 price + volume label update
Set-TextCell "D15" "0.0005978"
$ws.Range("E15").Value = "14OneONEWorstin24h"

Set-TextCell "D16" "0.006099"
Set-TextCell "D18" "2.092"
Set-TextCell "D19" "0.3198"
Set-TextCell "D20" "0.03142"
Set-TextCell "D22" "3.755"
Set-TextCell "D23" "0.04650"
Set-TextCell "D25" "0.001237"
Set-TextCell "D27" "0.00009700"
Set-TextCell "D28" "0.0001501"
Set-TextCell "D40" "0.03614"

# Rows 41-43: coin ordering rotated (41<-old43, 42<-old41, 43<-old42),
# prices independently refreshed
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextCell "D41" "0.006320"
$ws.Range("E41").Value = "40KickTokenKICK"

$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextCell "D42" "0.1048"
$ws.Range("E42").Value = "41BKEXTokenBKK"

$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextCell "D43" "0.002890"
$ws.Range("E43").Value = "42CEJICEJI"

Set-TextCell "D44" "0.007355"
Set-TextCell "D45" "0.00005282"

$ws.Range("E47").Value = "46CoinbaseStockTokenCOIN"

Set-TextCell "D48" "0.002311"
Set-TextCell "D49" "0.00002100"
Set-TextCell "D50" "0.0002000"
